# UserStory 5 Tasks erstellt
$wb = $excel.ActiveWorkbook

$tasks = $wb.Worksheets.Item("Tasks")
$prior = $wb.Worksheets.Item("Priorisierung")

# --- Add new rows 34-36 to the Tasks sheet ---

# Row 34: section header for UserStory 5 ("5." + the story text)
$tasks.Cells.Item(34, 1).Value = "5."
$tasks.Cells.Item(34, 2).Value = "Als Benutzer möchte ich den Standort des Hotels sehen und direkt zu Google Maps kommen"
$tasks.Rows.Item(34).RowHeight = 30

# Row 35 & 36: the new tasks
$tasks.Cells.Item(35, 2).Value = "Grundlegende Karte einbauen"
$tasks.Cells.Item(36, 2).Value = "Layout anpassen"

# --- Update the view/selection state on both sheets ---

$tasks.Activate()
$tasks.Application.ActiveWindow.ScrollRow = 22
$tasks.Rows.Item(36).EntireRow.Select()

$prior.Activate()
$prior.Range("B7").Select()

$tasks.Activate()
